$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (bold font, border, centered/top alignment) from A2
# into the new A4:A19 cells without disturbing their values (format-only paste).
$ws.Range("A2").Copy()
$ws.Range("A4:A19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "How do you approach making the most of study abroad opportunities for personal or professional growth through cultural exchange?"
$ws.Range("C4").Value = "I am going to pursue higher education in the US and am looking for suggestions for personal and professional growth. It's a part of cultural exchange. How can I make most from it ?"
$ws.Range("D4").Value = "2023-05-02 17:49:07.378971"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = "culture,higher studies,personal development"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "What kind of teamwork or collaboration skills have you developed through extracurricular activities or group projects?"
$ws.Range("C5").Value = "I want to join a samvad club of vnit.How can a club help to achieve soft skills"
$ws.Range("D5").Value = "2023-05-02 17:51:50.429037"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = "extracurricular,soft skills"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "What kind of extracurricular activities did you participate in while you were in college?"
$ws.Range("C6").Value = "There are several extracurricular activities that I want to take part in what should I choose?"
$ws.Range("D6").Value = "2023-05-02 19:13:29.922515"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = "campus life,extracurricular"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "How have you balanced pursuing higher education or certifications with working full-time, and what benefits have you seen from furthering your education in your career?"
$ws.Range("C7").Value = "How have you balanced pursuing higher education or certifications with working full-time, and what benefits have you seen from furthering your education in your career?"
$ws.Range("D7").Value = "2023-05-03 05:01:14.027499"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = "career,professional development,workplace culture,networking"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "What skills did you learn through your internships that have been most valuable in your current job, and how did you highlight those skills during the job application process?"
$ws.Range("C8").Value = "What skills did you learn through your internships that have been most valuable in your current job, and how did you highlight those skills during the job application process?"
$ws.Range("D8").Value = "2023-05-03 05:01:14.019426"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = "skills,job search/internship"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "How have you leveraged your network to identify job opportunities and gain insights into the job market, and how have those experiences helped you advance in your career?"
$ws.Range("C9").Value = "How have you leveraged your network to identify job opportunities and gain insights into the job market, and how have those experiences helped you advance in your career?"
$ws.Range("D9").Value = "2023-05-03 05:01:14.011040"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = "networking,job search/internship,career,professional development"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "How have your internships and other work experiences helped you develop transferable skills that are valuable across different industries, and how have you communicated those skills to potential employers during the job search process?"
$ws.Range("C10").Value = "How have your internships and other work experiences helped you develop transferable skills that are valuable across different industries, and how have you communicated those skills to potential employers during the job search process?"
$ws.Range("D10").Value = "2023-05-03 05:01:13.999975"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = "skills,job search/internship"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "How has your major prepared you for the job market, and what steps have you taken to develop transferable skills that are in high demand?"
$ws.Range("C11").Value = "How has your major prepared you for the job market, and what steps have you taken to develop transferable skills that are in high demand?"
$ws.Range("D11").Value = "2023-05-03 05:01:13.992465"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = "academics,career,skills"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "What challenges did you encounter during your job search, and how did you overcome them?"
$ws.Range("C12").Value = "What challenges did you encounter during your job search, and how did you overcome them?"
$ws.Range("D12").Value = "2023-05-03 05:01:13.985464"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = "job search/internship"

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "How have your internships and career preparation contributed to your ongoing career development?"
$ws.Range("C13").Value = "How have your internships and career preparation contributed to your ongoing career development?"
$ws.Range("D13").Value = "2023-05-03 05:01:13.978467"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = "career,job search/internship,career,professional development"

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "What strategies did you use to build a strong professional network during college, and how have those connections helped you in your career?"
$ws.Range("C14").Value = "What strategies did you use to build a strong professional network during college, and how have those connections helped you in your career?"
$ws.Range("D14").Value = "2023-05-03 05:01:13.972456"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = "networking,career,professional development"

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "How did you build professional connections and network with industry leaders during college, and how have those relationships helped you in your career?"
$ws.Range("C15").Value = "How did you build professional connections and network with industry leaders during college, and how have those relationships helped you in your career?"
$ws.Range("D15").Value = "2023-05-03 05:01:13.964700"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = "networking,career,professional development"

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "how did you make connections with fellow students and alumni during your time remaining in college?"
$ws.Range("C16").Value = "how did you make connections with fellow students and alumni during your time remaining in college?"
$ws.Range("D16").Value = "2023-05-03 05:01:13.951517"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = "networking,campus life"

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "How can I effectively balance coursework and part-time work during college?"
$ws.Range("C17").Value = "Is there any part-time work opportunity in college that might help me with pocket money?"
$ws.Range("D17").Value = "2023-05-03 05:13:05.906595"
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = "higher studies,time management"

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "How can I effectively balance coursework and part-time work during college?"
$ws.Range("C18").Value = "Is there any part-time work opportunity in college that might help me with pocket money?"
$ws.Range("D18").Value = "2023-05-03 05:13:30.888380"
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = "higher studies,time management"

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "How can I effectively prepare for and navigate a job interview or hiring process?"
$ws.Range("C19").Value = "How does T&P provide assistance in hiring process during College.."
$ws.Range("D19").Value = "2023-05-03 05:15:51.244709"
$ws.Range("E19").Value = 4
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = "career,industry skills"

